$wb = $excel.ActiveWorkbook

# --- Fix small spelling errors ---

# 1) "Department of Health, Western Australia (DOHWA)" -> "(DoHWA)"
#    in both the "experience" and "experience_latex" sheets (rows 9-14, 16, 17, column B).
$fixedInstitution = "Department of Health, Western Australia (DoHWA)"
$institutionRows = @(9, 10, 11, 12, 13, 14, 16, 17)

foreach ($sheetName in @("experience", "experience_latex")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $institutionRows) {
        $ws.Cells.Item($r, 2).Value = $fixedInstitution
    }
    # 2) "Course title: Applied Reggresssion Analyses (PUBH7631)" -> "Regression"
    $ws.Cells.Item(4, 5).Value = "Course title: Applied Regression Analyses (PUBH7631)"
}

# --- Active sheet / selection housekeeping (matches the author switching tabs) ---
$wsExperience = $wb.Worksheets.Item("experience")
$null = $wsExperience.Range("B20").Select()

$wsLatex = $wb.Worksheets.Item("experience_latex")
$null = $wsLatex.Activate()
$null = $wsLatex.Range("E4").Select()
